# Reuploaded fig2 and 3 in Developer guide
#
# 1) Update the two code-label text boxes on the sequence diagram slide.
# 2) Refresh the auto "last generated" date stamp (datetimeFigureOut field)
#    that lives on the slide master, every slide layout and the notes
#    master, from 12/15/2018 to 2/25/2019.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Rename deletePerson(p) -> deleteRecord(r) and
#    saveAddressBook(AddressBook) -> saveFinanceLog(FinanceLog)
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.HasText -eq -1) {
        $tr = $shp.TextFrame.TextRange
        $txt = $tr.Text

        if ($txt -eq "deletePerson(p)") {
            $tr.Replace("deletePerson", "deleteRecord") | Out-Null
            $tr.Replace("(p)", "(r)") | Out-Null
        }
        elseif ($txt -eq "saveAddressBook(AddressBook)") {
            $tr.Replace("saveAddressBook", "saveFinanceLog") | Out-Null
            $tr.Replace("AddressBook", "FinanceLog") | Out-Null
        }
    }
}

# ---------------------------------------------------------------------
# 2) Bump the cached date-field text wherever it appears: slide master,
#    every slide layout, and the notes master.
# ---------------------------------------------------------------------
$oldDate = "12/15/2018"
$newDate = "2/25/2019"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $isDate = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) { $isDate = $true }
            } catch { }

            if ($isDate -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $p.SlideMaster.CustomLayouts.Item($li)
}

Update-DatePlaceholder $p.NotesMaster
